$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 8190
$ws.Range("B4").Value = 7004
$ws.Range("B5").Value = 15194
$ws.Range("B6").Value = 0.0802775
